$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) - force text format to preserve exact string representation
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.276.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.617.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '195.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.212'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.20'
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.190.22'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '597.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.421.66'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.619.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.82'
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.27'
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.25'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0883'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.938.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '526.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0460'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.53'
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.56'
$ws.Range("D48").Style = "Normal"

# Update Volume(1h) column (E)
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  +4.21%  '
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("E18").Value = '  +3.07%  '
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("E20").Value = '  +1.65%  '
$ws.Range("E21").Value = '  +0.65%  '
$ws.Range("E22").Value = '  +1.50%  '
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("E24").Value = '  +0.50%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("E26").Value = '  -4.26%  '
$ws.Range("E27").Value = '  -2.96%  '
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("E29").Value = '  +0.93%  '
$ws.Range("E30").Value = '  +8.86%  '
$ws.Range("E31").Value = '  +2.75%  '
$ws.Range("E32").Value = '  -2.27%  '
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("E35").Value = '  +2.47%  '
$ws.Range("E36").Value = '  +5.29%  '
$ws.Range("E37").Value = '  +8.47%  '
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("E41").Value = '  -0.82%  '
$ws.Range("E42").Value = '  -2.27%  '
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("E45").Value = '  +6.49%  '
$ws.Range("E46").Value = '  +1.28%  '
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("E51").Value = '  +1.38%  '
